$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "WSUGSTIR-"

# Set B2 to an explicit empty text value (matching the other blank text
# cells in the row) rather than clearing it to a truly-blank cell.
$ws.Range("B2").Value = "'"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "29AAFCA0924K1ZN"
$ws.Range("H2").Value = 0.92
$ws.Range("I2").Value = 0.99
$ws.Range("J2").Value = 0.8
